# Auto-generated edit script
# Adds violent-crime data for 2025-05-16 (one additional day of data)
# by updating the 2025 (column L) and a few other year totals across
# the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 78
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 243

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 133
$ws.Range("L3").Value = 161
$ws.Range("L6").Value = 126
$ws.Range("L7").Value = 475

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 21
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 12
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L3").Value = 4
$ws.Range("L6").Value = 10

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L2").Value = 2
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 243
$ws.Range("L8").Value = 475
$ws.Range("L9").Value = 49
$ws.Range("L10").Value = 50
$ws.Range("L11").Value = 132
$ws.Range("L13").Value = 10
$ws.Range("L15").Value = 52
$ws.Range("L16").Value = 16
$ws.Range("L19").Value = 214
$ws.Range("L20").Value = 190
$ws.Range("L23").Value = 77
$ws.Range("L25").Value = 42
$ws.Range("L27").Value = 77
$ws.Range("L29").Value = 377
$ws.Range("L33").Value = 339
$ws.Range("K36").Value = 356
$ws.Range("L36").Value = 106
$ws.Range("K37").Value = 908
$ws.Range("L37").Value = 272
$ws.Range("L40").Value = 18
$ws.Range("L41").Value = 34
$ws.Range("L42").Value = 235
$ws.Range("L46").Value = 15
$ws.Range("L48").Value = 102
$ws.Range("L52").Value = 147
$ws.Range("L53").Value = 95
$ws.Range("L54").Value = 153
$ws.Range("L55").Value = 67
$ws.Range("E63").Value = 385
$ws.Range("J63").Value = 218
$ws.Range("K63").Value = 155
$ws.Range("L65").Value = 142
$ws.Range("L67").Value = 273
$ws.Range("L68").Value = 21
$ws.Range("L73").Value = 59
$ws.Range("L76").Value = 83
$ws.Range("L77").Value = 45
$ws.Range("L78").Value = 100
$ws.Range("L83").Value = 182
$ws.Range("L85").Value = 391
$ws.Range("J89").Value = 367
$ws.Range("L90").Value = 74
$ws.Range("L91").Value = 104
$ws.Range("L95").Value = 108
$ws.Range("L98").Value = 53
$ws.Range("L99").Value = 119
$ws.Range("E101").Value = 26056
$ws.Range("J101").Value = 29339
$ws.Range("K101").Value = 27556
$ws.Range("L101").Value = 7471

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 70
$ws.Range("L3").Value = 66
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 60
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 190

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2284
$ws.Range("J3").Value = 8079
$ws.Range("K3").Value = 8181
$ws.Range("L3").Value = 2322
$ws.Range("E4").Value = 2051
$ws.Range("K4").Value = 1764
$ws.Range("L4").Value = 626
$ws.Range("J5").Value = 630
$ws.Range("L5").Value = 140
$ws.Range("K6").Value = 9123
$ws.Range("L6").Value = 2099
$ws.Range("E7").Value = 26056
$ws.Range("J7").Value = 29339
$ws.Range("K7").Value = 27556
$ws.Range("L7").Value = 7471

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 30
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 135
$ws.Range("L7").Value = 377

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 89
$ws.Range("L6").Value = 118
$ws.Range("L7").Value = 339

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 114
$ws.Range("L3").Value = 27
$ws.Range("K7").Value = 356
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L5").Value = 11
$ws.Range("K6").Value = 272
$ws.Range("K7").Value = 908
$ws.Range("L7").Value = 272

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 71
$ws.Range("L7").Value = 235

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 33
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 80
$ws.Range("L7").Value = 153

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 78
$ws.Range("L3").Value = 93
$ws.Range("L5").Value = 7
$ws.Range("L7").Value = 273

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 74
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 391

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 103
$ws.Range("J7").Value = 367

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 29
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 119

Write-Output "Applied 2025-05-16 update across $($wb.Worksheets.Count) worksheet(s)."
